# Update column F (dSF) values per repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -10
    5  = -3
    6  = -1
    7  = -1
    8  = -3
    9  = 5
    13 = 0
    14 = 1
    17 = -5
    18 = 6
    21 = -1
    28 = 1
    32 = -9
    33 = 7
    34 = 14
    35 = -7
    36 = 13
    37 = -2
    39 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
